$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "1.00")
# are preserved exactly instead of being auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.554.87"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.642.95"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "212.88"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").Value = "0.532"
$ws.Range("E6").Value = "  +4.06%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "23.16"
$ws.Range("E8").Value = "  -5.30%  "
$ws.Range("D9").Value = "0.256"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Value = "0.0609"
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("D11").Value = "0.0888"
$ws.Range("E11").Value = "  +1.16%  "
$ws.Range("D12").Value = "1.877.42"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").Value = "1.657.37"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "0.580"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "4.01"
$ws.Range("E15").Value = "  -2.96%  "
$ws.Range("D16").Value = "64.19"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").Value = "27.590.82"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "229.10"
$ws.Range("E18").Value = "  -4.84%  "
$ws.Range("D19").Value = "0.0₃0722"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "7.52"
$ws.Range("E20").Value = "  -1.75%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  -3.44%  "
$ws.Range("D23").Value = "9.67"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("D25").Value = "149.18"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").Value = "6.96"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("D27").Value = "0.113"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "15.53"
$ws.Range("E29").Value = "  -4.64%  "
$ws.Range("D30").Value = "1.19"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Value = "0.0486"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").Value = "3.29"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").Value = "3.17"
$ws.Range("E33").Value = "  +2.31%  "
$ws.Range("D34").Value = "1.423.29"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  +1.07%  "
$ws.Range("D36").Value = "2.37"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "0.570"
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("D38").Value = "0.881"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("D39").Value = "0.0167"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").Value = "1.03"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "0.815"
$ws.Range("E42").Value = "  +3.20%  "
$ws.Range("D43").Value = "5.45"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "2.23"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "64.93"
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").Value = "1.786.21"
$ws.Range("E46").Value = "  -0.86%  "
$ws.Range("D47").Value = "1.66"
$ws.Range("E47").Value = "  -3.29%  "
$ws.Range("D48").Value = "86.98"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").Value = "0.0995"
$ws.Range("E50").Value = "  -2.63%  "
$ws.Range("D51").Value = "7.75"
$ws.Range("E51").Value = "  -1.30%  "

